# Replace the "Baz changes" / "Matthew Bortolin" paragraphs with a single
# new paragraph about version management.
$d = $word.ActiveDocument

# Locate the two paragraphs to be replaced by scanning for their text,
# rather than hard-coding paragraph indices.
$bazIndex = 0
$bortolinIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Baz changes*" -and $bazIndex -eq 0) {
        $bazIndex = $i
    }
    if ($txt -like "Matthew Bortolin*" -and $bortolinIndex -eq 0) {
        $bortolinIndex = $i
    }
}

$pBaz = $d.Paragraphs.Item($bazIndex)
$pBortolin = $d.Paragraphs.Item($bortolinIndex)

# Remove both paragraphs (including their paragraph marks) entirely.
$killRange = $d.Range($pBaz.Range.Start, $pBortolin.Range.End)
$killRange.Delete()

# Insert a fresh paragraph mark in their place so the document keeps the
# same total paragraph count (one new paragraph replacing the two removed
# ones), then fill it in with the replacement content/runs below.
$newPara = $d.Paragraphs.Item($bazIndex)
$newPara.Range.InsertParagraphBefore()
$target = $d.Paragraphs.Item($bazIndex)

$apos = [char]0x2019
$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B56A30" w:rsidRDefault="00B56A30">' +
    '<w:r><w:t xml:space="preserve">There will be many versions of software during development and only some are released to the users. </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Version management involves keeping track of multiple versions of the components and ensures that changes made </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>by  different</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> developers don' + $apos + 't interfere with eachother.</w:t></w:r>' +
    '</w:p>'

$target.Range.InsertXML($bodyXml) | Out-Null
